# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-03 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-04 Thursday", 2) | Out-Null

# The worksheet table's first data row effectively drops the leading
# "69÷6=" problem (its cell is removed) and gains a new "18÷3=" problem
# appended at the end; the three problems in between are also replaced
# with new values. Because every cell in the row shares identical
# formatting, this nets out to simply shifting the five cell values one
# slot to the left and filling in the new set of problems left-to-right.
$t = $d.Tables.Item(1)

$rowUpdates = @{
    1  = @("24÷2=", "52÷6=", "49÷8=", "14÷7=", "18÷3=")
    5  = @("25÷7=", "51÷5=", "16÷4=", "91÷4=", "72÷2=")
    9  = @("94÷3=", "64÷5=", "80÷3=", "41÷8=", "82÷7=")
    13 = @("78÷8=", "77÷7=", "34÷8=", "18÷2=", "80÷2=")
    17 = @("15÷9=", "80÷8=", "10÷4=", "68÷2=", "78÷4=")
}

foreach ($rowIndex in $rowUpdates.Keys) {
    $values = $rowUpdates[$rowIndex]
    $row = $t.Rows.Item($rowIndex)
    for ($col = 1; $col -le $values.Count; $col++) {
        $row.Cells.Item($col).Range.Text = $values[$col - 1]
    }
}
